$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns for the refreshed crypto data.
# D-column values are written with a leading apostrophe and the style is then
# reset to "Normal" so Excel keeps them as plain text (matching the original
# inlineStr cells) instead of auto-converting numeric-looking text into numbers.

$cell = $ws.Range("D2")
$cell.Value = '''62.327.56'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -2.92%  '

$cell = $ws.Range("D3")
$cell.Value = '''3.006.46'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -3.22%  '

$cell = $ws.Range("D4")
$cell.Value = '''0.999'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

$cell = $ws.Range("D5")
$cell.Value = '''581.52'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.80%  '

$cell = $ws.Range("D6")
$cell.Value = '''146.94'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -6.61%  '

$ws.Range("E7").Value = '  +0.16%  '

$cell = $ws.Range("D8")
$cell.Value = '''0.523'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.40%  '

$cell = $ws.Range("D9")
$cell.Value = '''3.004.42'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -3.07%  '

$cell = $ws.Range("D10")
$cell.Value = '''0.149'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -6.51%  '

$cell = $ws.Range("D11")
$cell.Value = '''5.66'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -4.57%  '

$cell = $ws.Range("D12")
$cell.Value = '''0.442'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -2.49%  '

$cell = $ws.Range("D13")
$cell.Value = '''0.0000228'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -5.15%  '

$cell = $ws.Range("D14")
$cell.Value = '''34.64'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -6.93%  '

$cell = $ws.Range("D15")
$cell.Value = '''0.122'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.78%  '

$cell = $ws.Range("D16")
$cell.Value = '''3.497.34'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -2.84%  '

$cell = $ws.Range("D17")
$cell.Value = '''7.08'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -2.21%  '

$cell = $ws.Range("D18")
$cell.Value = '''62.300.31'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.49%  '

$cell = $ws.Range("D19")
$cell.Value = '''2.997.07'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -2.85%  '

$cell = $ws.Range("D20")
$cell.Value = '''453.45'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -4.51%  '

$cell = $ws.Range("D21")
$cell.Value = '''13.87'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -3.98%  '

$cell = $ws.Range("D22")
$cell.Value = '''0.680'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -4.38%  '

$cell = $ws.Range("D23")
$cell.Value = '''7.31'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -3.64%  '

$cell = $ws.Range("D24")
$cell.Value = '''2.30'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -6.16%  '

$cell = $ws.Range("D25")
$cell.Value = '''80.14'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -1.01%  '

$cell = $ws.Range("D26")
$cell.Value = '''12.33'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -4.87%  '

$cell = $ws.Range("D27")
$cell.Value = '''10.10'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -2.93%  '

$cell = $ws.Range("D28")
$cell.Value = '''1.00'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.02%  '

$cell = $ws.Range("D30")
$cell.Value = '''7.15'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -4.40%  '

$cell = $ws.Range("D31")
$cell.Value = '''2.61'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.74%  '

$cell = $ws.Range("D32")
$cell.Value = '''2.09'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -4.21%  '

$cell = $ws.Range("D33")
$cell.Value = '''26.86'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.01%  '

$cell = $ws.Range("D34")
$cell.Value = '''0.107'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -4.98%  '

$cell = $ws.Range("D36")
$cell.Value = '''0.0₃0792'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -6.25%  '

$cell = $ws.Range("D37")
$cell.Value = '''5.74'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -4.97%  '

$cell = $ws.Range("D38")
$cell.Value = '''2.13'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -5.15%  '

$cell = $ws.Range("D39")
$cell.Value = '''50.16'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -1.48%  '

$cell = $ws.Range("D40")
$cell.Value = '''9.04'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '

$cell = $ws.Range("D41")
$cell.Value = '''2.94'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -11.05%  '

$cell = $ws.Range("D42")
$cell.Value = '''407.44'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -7.16%  '

$cell = $ws.Range("D43")
$cell.Value = '''0.277'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -4.84%  '

$ws.Range("E44").Value = '  -0.84%  '

$cell = $ws.Range("D45")
$cell.Value = '''2.767.78'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.01%  '

$cell = $ws.Range("D46")
$cell.Value = '''0.0352'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -3.27%  '

$cell = $ws.Range("D47")
$cell.Value = '''38.03'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -5.81%  '

$cell = $ws.Range("D48")
$cell.Value = '''127.87'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -1.72%  '

$ws.Range("E49").Value = '  +0.03%  '

$cell = $ws.Range("D50")
$cell.Value = '''0.108'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -1.88%  '

$cell = $ws.Range("D51")
$cell.Value = '''23.84'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -6.06%  '
